$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared string "ECs" will be introduced naturally via the string values below.

# Row 2: FAPs / Nrg2 / Erbb3 / ECs
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Nrg2"
$ws.Cells.Item(2, 3).Value = "Erbb3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.178245666666667
$ws.Cells.Item(2, 8).Value = 3.534737
$ws.Cells.Item(2, 9).Value = 0.8546990546349293
$ws.Cells.Item(2, 10).Value = 0.8546990546349293
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.058393
$ws.Cells.Item(2, 14).Value = 0.175179
$ws.Cells.Item(2, 15).Value = 0.01085373024912483
$ws.Cells.Item(2, 16).Value = 0.01085373024912483
$ws.Cells.Item(2, 17).Value = 0.06880129921366666
$ws.Cells.Item(2, 18).Value = 0.6192116929229999
$ws.Cells.Item(2, 19).Value = 0.009276672983189528
$ws.Cells.Item(2, 20).Value = 0.009276672983189528

# Row 3: FAPs / Nrg2 / Erbb3 / FAPs
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Nrg2"
$ws.Cells.Item(3, 3).Value = "Erbb3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.178245666666667
$ws.Cells.Item(3, 8).Value = 3.534737
$ws.Cells.Item(3, 9).Value = 0.8546990546349293
$ws.Cells.Item(3, 10).Value = 0.8546990546349293
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.4773683333333333
$ws.Cells.Item(3, 14).Value = 1.432105
$ws.Cells.Item(3, 15).Value = 0.08873027793527143
$ws.Cells.Item(3, 16).Value = 0.08873027793527143
$ws.Cells.Item(3, 17).Value = 0.5624571701538889
$ws.Cells.Item(3, 18).Value = 5.062114531384999
$ws.Cells.Item(3, 19).Value = 0.07583768466877101
$ws.Cells.Item(3, 20).Value = 0.07583768466877101

# Row 4: FAPs / Nrg2 / Erbb3 / sCs
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Nrg2"
$ws.Cells.Item(4, 3).Value = "Erbb3"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.178245666666667
$ws.Cells.Item(4, 8).Value = 3.534737
$ws.Cells.Item(4, 9).Value = 0.8546990546349293
$ws.Cells.Item(4, 10).Value = 0.8546990546349293
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 4.844232333333333
$ws.Cells.Item(4, 14).Value = 14.532697
$ws.Cells.Item(4, 15).Value = 0.9004159918156037
$ws.Cells.Item(4, 16).Value = 0.9004159918156038
$ws.Cells.Item(4, 17).Value = 5.707695755076554
$ws.Cells.Item(4, 18).Value = 51.369261795689
$ws.Cells.Item(4, 19).Value = 0.7695846969829688
$ws.Cells.Item(4, 20).Value = 0.7695846969829688

# Row 5: sCs / Nrg2 / Erbb3 / ECs
$ws.Cells.Item(5, 1).Value = "sCs"
$ws.Cells.Item(5, 2).Value = "Nrg2"
$ws.Cells.Item(5, 3).Value = "Erbb3"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.2003046666666667
$ws.Cells.Item(5, 8).Value = 0.6009139999999999
$ws.Cells.Item(5, 9).Value = 0.1453009453650707
$ws.Cells.Item(5, 10).Value = 0.1453009453650707
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.058393
$ws.Cells.Item(5, 14).Value = 0.175179
$ws.Cells.Item(5, 15).Value = 0.01085373024912483
$ws.Cells.Item(5, 16).Value = 0.01085373024912483
$ws.Cells.Item(5, 17).Value = 0.01169639040066667
$ws.Cells.Item(5, 18).Value = 0.105267513606
$ws.Cells.Item(5, 19).Value = 0.001577057265935302
$ws.Cells.Item(5, 20).Value = 0.001577057265935302

# Row 6: sCs / Nrg2 / Erbb3 / FAPs
$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Nrg2"
$ws.Cells.Item(6, 3).Value = "Erbb3"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.2003046666666667
$ws.Cells.Item(6, 8).Value = 0.6009139999999999
$ws.Cells.Item(6, 9).Value = 0.1453009453650707
$ws.Cells.Item(6, 10).Value = 0.1453009453650707
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.4773683333333333
$ws.Cells.Item(6, 14).Value = 1.432105
$ws.Cells.Item(6, 15).Value = 0.08873027793527143
$ws.Cells.Item(6, 16).Value = 0.08873027793527143
$ws.Cells.Item(6, 17).Value = 0.09561910488555556
$ws.Cells.Item(6, 18).Value = 0.8605719439699999
$ws.Cells.Item(6, 19).Value = 0.01289259326650041
$ws.Cells.Item(6, 20).Value = 0.01289259326650041

# Row 7: sCs / Nrg2 / Erbb3 / sCs
$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Nrg2"
$ws.Cells.Item(7, 3).Value = "Erbb3"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.2003046666666667
$ws.Cells.Item(7, 8).Value = 0.6009139999999999
$ws.Cells.Item(7, 9).Value = 0.1453009453650707
$ws.Cells.Item(7, 10).Value = 0.1453009453650707
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 4.844232333333333
$ws.Cells.Item(7, 14).Value = 14.532697
$ws.Cells.Item(7, 15).Value = 0.9004159918156037
$ws.Cells.Item(7, 16).Value = 0.9004159918156038
$ws.Cells.Item(7, 17).Value = 0.9703223427842221
$ws.Cells.Item(7, 18).Value = 8.732901085058
$ws.Cells.Item(7, 19).Value = 0.130831294832635
$ws.Cells.Item(7, 20).Value = 0.130831294832635
